$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.883.21'
$ws.Range("E2").Value = '  -0.24%  '

$ws.Range("D3").Value = '2.355.62'
$ws.Range("E3").Value = '  -1.96%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").Value = '''318.40'
$ws.Range("E5").Value = '  -4.77%  '

$ws.Range("D6").Value = '''107.04'
$ws.Range("E6").Value = '  +4.42%  '

$ws.Range("E7").Value = '  -1.58%  '

$ws.Range("E8").Value = '  +0.12%  '

$ws.Range("D9").Value = '''0.620'
$ws.Range("E9").Value = '  -3.18%  '

$ws.Range("D10").Value = '''41.24'

$ws.Range("D11").Value = '''0.0927'
$ws.Range("E11").Value = '  -0.95%  '

$ws.Range("D12").Value = '''8.53'
$ws.Range("E12").Value = '  -0.71%  '

$ws.Range("D13").Value = '''0.999'
$ws.Range("E13").Value = '  -4.61%  '

$ws.Range("E14").Value = '  -0.24%  '

$ws.Range("D15").Value = '''15.97'
$ws.Range("E15").Value = '  -5.47%  '

$ws.Range("D16").Value = '2.713.83'
$ws.Range("E16").Value = '  -1.75%  '

$ws.Range("D17").Value = '2.322.21'
$ws.Range("E17").Value = '  -2.91%  '

$ws.Range("D18").Value = '42.821.22'
$ws.Range("E18").Value = '  -0.27%  '

$ws.Range("D19").Value = '''7.56'
$ws.Range("E19").Value = '  -0.23%  '

$ws.Range("D20").Value = '''0.0000107'
$ws.Range("E20").Value = '  -1.28%  '

$ws.Range("D21").Value = '''76.08'
$ws.Range("E21").Value = '  -0.74%  '

$ws.Range("D22").Value = '''3.63'
$ws.Range("E22").Value = '  -7.68%  '

$ws.Range("D23").Value = '''268.62'
$ws.Range("E23").Value = '  -0.92%  '

$ws.Range("D24").Value = '''2.31'
$ws.Range("E24").Value = '  -3.53%  '

$ws.Range("D25").Value = '''9.41'
$ws.Range("E25").Value = '  -7.30%  '

$ws.Range("D26").Value = '''0.999'
$ws.Range("E26").Value = '  -0.14%  '

$ws.Range("E27").Value = '  -3.67%  '

$ws.Range("D28").Value = '''23.47'
$ws.Range("E28").Value = '  -2.54%  '

$ws.Range("E29").Value = '  +2.12%  '

$ws.Range("D30").Value = '''36.86'
$ws.Range("E30").Value = '  +1.09%  '

$ws.Range("D31").Value = '''167.50'
$ws.Range("E31").Value = '  -3.76%  '

$ws.Range("D32").Value = '''0.0907'
$ws.Range("E32").Value = '  -1.94%  '

$ws.Range("E33").Value = '  -0.67%  '

$ws.Range("D34").Value = '''2.89'
$ws.Range("E34").Value = '  -7.01%  '

$ws.Range("E35").Value = '  +11.98%  '

$ws.Range("E36").Value = '  -2.44%  '

$ws.Range("D37").Value = '''4.74'
$ws.Range("E37").Value = '  -1.60%  '

$ws.Range("E38").Value = '  -0.48%  '

$ws.Range("D39").Value = '''3.82'
$ws.Range("E39").Value = '  -4.01%  '

$ws.Range("D40").Value = '''2.72'
$ws.Range("E40").Value = '  -5.89%  '

$ws.Range("D41").Value = '''106.81'
$ws.Range("E41").Value = '  +14.31%  '

$ws.Range("E42").Value = '  -2.80%  '

$ws.Range("D43").Value = '''0.239'
$ws.Range("E43").Value = '  +2.05%  '

$ws.Range("E44").Value = '  +1.17%  '

$ws.Range("E45").Value = '  +0.06%  '

$ws.Range("D46").Value = '''12.32'
$ws.Range("E46").Value = '  +1.66%  '

$ws.Range("D47").Value = '''113.39'
$ws.Range("E47").Value = '  -4.57%  '

$ws.Range("D48").Value = '''5.52'
$ws.Range("E48").Value = '  -0.38%  '

$ws.Range("D49").Value = '''9.15'
$ws.Range("E49").Value = '  -0.12%  '

$ws.Range("D50").Value = '''75.64'
$ws.Range("E50").Value = '  +7.12%  '

$ws.Range("D51").Value = '''1.29'
$ws.Range("E51").Value = '  -0.38%  '
